# Generate Report for Handoff
# Source file identifier changed from
#   de258464-8020-4224-862f-d05d1444577b
# to
#   0fd03137-d045-4f89-b608-541ceb656119
# and associated handoff timestamps/xliff hashes were refreshed.

$wb = $excel.ActiveWorkbook

$oldGuid = "de258464-8020-4224-862f-d05d1444577b"
$newGuid = "0fd03137-d045-4f89-b608-541ceb656119"

$newZhXlf = "$newGuid.13a06edb5b6890ad6752081af87ce7369f7f7275.zh-cn.xlf"
$newDeXlf = "$newGuid.13a06edb5b6890ad6752081af87ce7369f7f7275.de-de.xlf"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/oltest/blob/d8aa3e536b6be30dfd2095cb965776c1835845dc/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name / Path And Name / Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-07-26 08:07:51"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name / Latest Handoff File / Latest Handoff Datetime
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("F2").Value = $newZhXlf
$wsZh.Range("G2").Value = "2016-07-26 08:07:41"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

# ---------------------------------------------------------------------------
# Sheet "de-de": Source File Name / Latest Handoff File / Latest Handoff Datetime
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("F2").Value = $newDeXlf
$wsDe.Range("G2").Value = "2016-07-26 08:07:51"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")
